# Applies the cell-value updates for the Seraph_Profits scheduled-runner sync.
# All sheets hold static (non-formula) cached market-board figures; each hunk
# below corresponds to one Leve row whose NQ/HQ price/profit columns (H:N) were
# refreshed. Empty cells in the source data are genuinely absent (no formula),
# so ClearContents() is used wherever a column disappears in the diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 160
$ws.Range("I33").Value = 160
$ws.Range("K33").Value = 160
$ws.Range("M33").Value = 69

$ws.Range("H41").Value = 222
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 222
$ws.Range("K41").Value = 0
$ws.Range("L41").ClearContents()
$ws.Range("M41").Value = 222
$ws.Range("N41").Value = -1102

$ws.Range("I51").Value = 11285.714
$ws.Range("J51").Value = 12278.8
$ws.Range("K51").Value = 11285.714
$ws.Range("L51").Value = 12278.8
$ws.Range("M51").Value = -10801.714
$ws.Range("N51").Value = -13246.8

$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

$ws.Range("H76").Value = 7500.5
$ws.Range("I76").Value = 5003
$ws.Range("K76").Value = 5003
$ws.Range("M76").Value = -4688

$ws.Range("H79").Value = 7500.5
$ws.Range("I79").Value = 5003
$ws.Range("K79").Value = 5003
$ws.Range("M79").Value = -3911

$ws.Range("H80").Value = 13414.833
$ws.Range("J80").Value = 13414.833
$ws.Range("L80").Value = 40244.499
$ws.Range("N80").Value = -42240.499

$ws.Range("H83").Value = 13414.833
$ws.Range("J83").Value = 13414.833
$ws.Range("L83").Value = 120733.497
$ws.Range("N83").Value = -130717.497

$ws.Range("H100").Value = 1062.6
$ws.Range("I100").Value = 1069.5555
$ws.Range("K100").Value = 1069.5555
$ws.Range("M100").Value = -528.5554999999999

$ws.Range("H106").Value = 43496204
$ws.Range("I106").Value = 58839864
$ws.Range("J106").Value = 22499.5
$ws.Range("K106").Value = 58839864
$ws.Range("L106").Value = 22499.5
$ws.Range("M106").Value = -58839233
$ws.Range("N106").Value = -23761.5

$ws.Range("H132").Value = 3379.3572
$ws.Range("I132").Value = 2704.8262
$ws.Range("K132").Value = 8114.4786
$ws.Range("M132").Value = -5584.4786

$ws.Range("H137").Value = 3089.9048
$ws.Range("I137").Value = 3164
$ws.Range("K137").Value = 9492
$ws.Range("M137").Value = -6942


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 166667170
$ws.Range("I2").Value = 333333340
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 333333340
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -333333227
$ws.Range("N2").Value = -1226

$ws.Range("H6").Value = 97263.45
$ws.Range("I6").Value = 133249.88
$ws.Range("J6").Value = 1299.6666
$ws.Range("K6").Value = 133249.88
$ws.Range("L6").Value = 1299.6666
$ws.Range("M6").Value = -133076.88
$ws.Range("N6").Value = -1645.6666

$ws.Range("H12").Value = 252724.83
$ws.Range("J12").Value = 3812.375
$ws.Range("L12").Value = 3812.375
$ws.Range("N12").Value = -4158.375

$ws.Range("H74").Value = 1505.9656
$ws.Range("I74").Value = 1426.24
$ws.Range("J74").Value = 2004.25
$ws.Range("K74").Value = 1426.24
$ws.Range("L74").Value = 2004.25
$ws.Range("M74").Value = -552.24
$ws.Range("N74").Value = -3752.25

$ws.Range("H77").Value = 1505.9656
$ws.Range("I77").Value = 1426.24
$ws.Range("J77").Value = 2004.25
$ws.Range("K77").Value = 7131.2
$ws.Range("L77").Value = 10021.25
$ws.Range("M77").Value = -2763.2
$ws.Range("N77").Value = -18757.25

$ws.Range("H116").Value = 166667170
$ws.Range("I116").Value = 333333340
$ws.Range("J116").Value = 1000
$ws.Range("K116").Value = 333333340
$ws.Range("L116").Value = 1000
$ws.Range("M116").Value = -333331046
$ws.Range("N116").Value = -5588


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 166667170
$ws.Range("I3").Value = 333333340
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 333333340
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -333333226
$ws.Range("N3").Value = -1228

$ws.Range("H7").Value = 364.125
$ws.Range("I7").Value = 271.77777
$ws.Range("J7").Value = 482.85715
$ws.Range("K7").Value = 271.77777
$ws.Range("L7").Value = 482.85715
$ws.Range("M7").Value = -158.77777
$ws.Range("N7").Value = -708.85715

$ws.Range("H11").Value = 499.16666
$ws.Range("J11").Value = 685.625
$ws.Range("L11").Value = 685.625
$ws.Range("N11").Value = -965.625

$ws.Range("H25").Value = 3302.1
$ws.Range("I25").Value = 3833.5
$ws.Range("J25").Value = 2505
$ws.Range("K25").Value = 3833.5
$ws.Range("L25").Value = 2505
$ws.Range("M25").Value = -3598.5
$ws.Range("N25").Value = -2975

$ws.Range("H36").Value = 1423.125
$ws.Range("I36").Value = 1423.125
$ws.Range("K36").Value = 1423.125
$ws.Range("M36").Value = -889.125

$ws.Range("H54").Value = 16666.666
$ws.Range("I54").Value = 22500
$ws.Range("J54").Value = 5000
$ws.Range("K54").Value = 22500
$ws.Range("L54").Value = 5000
$ws.Range("M54").Value = -22016
$ws.Range("N54").Value = -5968

$ws.Range("H134").Value = 2900.7144
$ws.Range("I134").Value = 2900.7144
$ws.Range("K134").Value = 8702.143199999999
$ws.Range("M134").Value = -6167.143199999999


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 10600.4

$ws.Range("H5").Value = 1415
$ws.Range("I5").Value = 2025
$ws.Range("K5").Value = 2025
$ws.Range("M5").Value = -1913

$ws.Range("H31").Value = 3133.6758
$ws.Range("I31").Value = 3423.75
$ws.Range("K31").Value = 3423.75
$ws.Range("M31").Value = -3128.75

$ws.Range("H34").Value = 3133.6758
$ws.Range("I34").Value = 3423.75
$ws.Range("K34").Value = 3423.75
$ws.Range("M34").Value = -3221.75

$ws.Range("H47").Value = 4000
$ws.Range("I47").Value = 4000
$ws.Range("K47").Value = 4000
$ws.Range("M47").Value = -3434

$ws.Range("H103").Value = 47450
$ws.Range("I103").Value = 34900
$ws.Range("K103").Value = 34900
$ws.Range("M103").Value = -33728

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

$ws.Range("H132").Value = 4357.4287
$ws.Range("I132").Value = 4175.3335
$ws.Range("K132").Value = 12526.0005
$ws.Range("M132").Value = -9996.000499999998

$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6166724.5
$ws.Range("I4").Value = 7536708.5
$ws.Range("J4").Value = 138793.2
$ws.Range("K4").Value = 22610125.5
$ws.Range("L4").Value = 416379.6
$ws.Range("M4").Value = -22610013.5
$ws.Range("N4").Value = -416603.6

$ws.Range("H64").Value = 83343944
$ws.Range("I64").Value = 1220
$ws.Range("K64").Value = 3660
$ws.Range("M64").Value = -3390

$ws.Range("H67").Value = 83343944
$ws.Range("I67").Value = 1220
$ws.Range("K67").Value = 3660
$ws.Range("M67").Value = -2724

$ws.Range("H86").Value = 1649.2
$ws.Range("I86").Value = 1415.6666
$ws.Range("K86").Value = 4246.9998
$ws.Range("M86").Value = -3060.9998

$ws.Range("H89").Value = 1649.2
$ws.Range("I89").Value = 1415.6666
$ws.Range("K89").Value = 12740.9994
$ws.Range("M89").Value = -6812.999400000001

$ws.Range("H98").Value = 2005.6
$ws.Range("I98").Value = 2005.6
$ws.Range("K98").Value = 6016.799999999999
$ws.Range("M98").Value = -4518.799999999999

$ws.Range("H113").Value = 908.375
$ws.Range("J113").Value = 940.8
$ws.Range("L113").Value = 2822.4
$ws.Range("N113").Value = -7162.4

$ws.Range("H132").Value = 2499.2222
$ws.Range("J132").Value = 2701
$ws.Range("L132").Value = 24309
$ws.Range("N132").Value = -29369


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3076.4707
$ws.Range("I132").Value = 2275
$ws.Range("K132").Value = 6825
$ws.Range("M132").Value = -4295


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3606400
$ws.Range("I2").Value = 4000800
$ws.Range("J2").Value = 3211999.8
$ws.Range("K2").Value = 4000800
$ws.Range("L2").Value = 3211999.8
$ws.Range("M2").Value = -4000688
$ws.Range("N2").Value = -3212223.8

$ws.Range("H45").Value = 37999.332
$ws.Range("J45").Value = 46999
$ws.Range("L45").Value = 46999
$ws.Range("N45").Value = -47813

$ws.Range("H48").Value = 9000
$ws.Range("I48").Value = 9000
$ws.Range("K48").Value = 9000
$ws.Range("M48").Value = -8339

$ws.Range("H56").Value = 61683
$ws.Range("I56").Value = 61683
$ws.Range("K56").Value = 61683
$ws.Range("M56").Value = -60992


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 4601.4443
$ws.Range("J23").Value = 2714.75
$ws.Range("L23").Value = 2714.75
$ws.Range("N23").Value = -3172.75

$ws.Range("H96").Value = 3513.7334
$ws.Range("I96").Value = 3358.5
$ws.Range("K96").Value = 3358.5
$ws.Range("M96").Value = -1985.5

$ws.Range("H107").Value = 2416.3333
$ws.Range("I107").Value = 2249.75
$ws.Range("K107").Value = 6749.25
$ws.Range("M107").Value = -4829.25

$ws.Range("H140").Value = 57500
$ws.Range("J140").Value = 57500
$ws.Range("L140").Value = 57500
$ws.Range("N140").Value = -67860

